# Update countries & provincias Spain
#
# This updates the COVID stats table on the "Pais" sheet with a newer
# data snapshot. Totals were refreshed (column B "Casos totales" and
# dependents), and because the sheet is sorted by total cases
# descending, a handful of countries swap ranks/rows as a result:
#   - Banglades overtakes China (rows 21/22)
#   - Polonia overtakes Argentina (rows 39/40)
#   - Filipinas overtakes Irlanda (rows 41/42)
#   - Oman overtakes Republica Dominicana, which overtakes Rumania (rows 44/45/46)
# A few other rows (63, 103, 105) just get refreshed figures without a
# rank change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($Row, $Country, $Total, $New, $Active, $Recovered, $Critical, $DeathsToday, $Deaths)

    if ($Country -ne $null) {
        $ws.Cells.Item($Row, 1).Value = $Country
    }
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $New
    $ws.Cells.Item($Row, 4).Value = $Active
    $ws.Cells.Item($Row, 5).Value = $Recovered
    $ws.Cells.Item($Row, 6).Value = $Critical
    $ws.Cells.Item($Row, 7).Value = $DeathsToday
    $ws.Cells.Item($Row, 8).Value = $Deaths
}

# Banglades / China swap rank (rows 21-22)
Set-Row 21 "Banglades" 84379 2856 17827 65413 0 44 1139
Set-Row 22 "China"     83075 11   78367 74    0 0  4634

# Polonia / Argentina swap rank (rows 39-40)
Set-Row 39 "Polonia"   29017 440 14104 13676 0 15 1237
Set-Row 40 "Argentina" 28764 0   8743  19236 0 0  785

# Filipinas / Irlanda swap rank (rows 41-42)
Set-Row 41 "Filipinas" 25392 605 5706  18612 0 22 1074
Set-Row 42 "Irlanda"   25250 0   22698 847   0 0  1705

# Row 43 (Afganistan) keeps its rank and values unchanged.

# Oman / Republica Dominicana / Rumania shift rank (rows 44-46)
Set-Row 44 "Oman"                  22077 1006 7530  14448 0 3 99
Set-Row 45 "Republica Dominicana"  22008 0    12754 8686  0 0 568
Set-Row 46 "Rumania"               21404 0    15445 4579  0 0 1380

# Data-only refreshes (no rank change)
Set-Row 63  $null 9942 4 7215 2398 0 0 329
Set-Row 103 $null 1880 0 1252 617  0 0 11
Set-Row 105 $null 1763 7 1416 272  0 1 75
